# Added replacement feet for 1390-B
# - Insert a new line item row (new row 73) in the HARDWARE RNGs table for
#   "Replacement feet for 1390-B" (McMaster, 2 @ $5.96 = $11.92).
# - Update the price of the "Type 1390-B Random Noise Generator" line (row 72)
#   from 59.99 to 49.
# - Move the view/selection to the edited area.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the 1390-B generator unit price (row 72, still row 72 before insert)
$ws.Range("B72").Value = 49

# 2) Insert a brand-new row right after it (becomes row 73); this shifts every
#    row below (the subtotal row and all following tables) down by one and
#    keeps all formula references in sync, just like Excel's own Insert.
$ws.Rows.Item(73).Insert()

# The insert copies row 72's formatting into every column of the new row,
# including the (unwanted) BOUGHT?/ARRIVED? cells in E73:F73 that this line
# item doesn't use - clear those so the row only has the cells it needs.
$ws.Range("E73:F73").Clear()

# 3) Populate the new row with the replacement-feet line item.
$ws.Range("A73").Value = "Replacement feet for 1390-B"
$ws.Range("B73").Value = 5.96
$ws.Range("C73").Value = 2
$ws.Range("D73").Formula = "=C73*B73"
$ws.Range("G73").Value = "McMaster"
$ws.Range("I73").Value = "http://www.mcmaster.com/#9546k549/=y4zioh"

# 4) Update the view so the newly added row is visible/selected.
$excel.ActiveWindow.ScrollRow = 60
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A73").Select()
